$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time formatting used by A2:A3 down onto the new rows A4:A5
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New readxl example rows
$ws.Range("A4").Value = 42737.479166666664
$ws.Range("A5").Value = 42738.479166666664

# Shared formula across B4:B5, mirroring the existing TEXT() formula in B2:B3
$ws.Range("B4:B5").Formula = '=TEXT(A4,"yyyy-mm-dd hh:mm:ss")'

# Update the active selection as recorded in the diff
$ws.Range("A2").Select() | Out-Null
